$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$app = $ws.Application
$app.Goto($ws.Range("F9"), $true)
$win = $app.ActiveWindow
Write-Host "After ScrollRow:" $win.ScrollRow() "ScrollColumn:" $win.ScrollColumn()
